# Raul's Log.xlsx -- re-worked school zoning/distribution: appends new
# log entries (rows 701-709) to the "Logs" sheet, mirroring the existing
# row layout/styles (columns A-F) and updates the view/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# Each entry: row number, A..F text/values (F = $null means "leave blank")
$entries = @(
    @{ Row = 701; A = "Setup Mic";           B = 42656; C = "1800"; D = "DB";  E = "2027"; F = "Neck mic and small PA from DB 0003" },
    @{ Row = 702; A = "Pickup Mic";          B = 42656; C = "2100"; D = "DB";  E = "2027"; F = "Return neck mic and small PA to DB 0003" },
    @{ Row = 703; A = "Demo";                B = 42656; C = "1600"; D = "KT";  E = "519";  F = $null },
    @{ Row = 704; A = "Demo";                B = 42656; C = "1900"; D = "SSB"; E = "N108"; F = "Demo neck mic" },
    @{ Row = 705; A = "Other";               B = 42656; C = "1550"; D = "MC";  E = "101A"; F = "Please pick up wireless keyboard and remote for the projector, return it to FC 164" },
    @{ Row = 706; A = "SCLD Student Event";  B = 42656; C = "1800"; D = "ACW"; E = "209";  F = "INC000000733578" },
    @{ Row = 707; A = "SCLD Student Logout"; B = 42656; C = "2100"; D = "ACW"; E = "209";  F = "INC000000733578" },
    @{ Row = 708; A = "SCLD Student Event";  B = 42656; C = "1800"; D = "CLH"; E = "K";    F = "INC000000733917" },
    @{ Row = 709; A = "SCLD Student Logout"; B = 42656; C = "2030"; D = "CLH"; E = "K";    F = "INC000000733917" }
)

foreach ($entry in $entries) {
    $r = $entry.Row
    $ws.Cells.Item($r, 1).Value = $entry.A
    $ws.Cells.Item($r, 2).Value = $entry.B
    $ws.Cells.Item($r, 3).Value = $entry.C
    $ws.Cells.Item($r, 4).Value = $entry.D
    $ws.Cells.Item($r, 5).Value = $entry.E
    if ($entry.F -ne $null) {
        $ws.Cells.Item($r, 6).Value = $entry.F
    }
}

# Row 705 wraps onto two lines (long note in column F), matching the
# taller row height used elsewhere in the sheet for similar notes.
$ws.Rows.Item(705).RowHeight = 30

# Move the view/selection the way the author left it after the edit.
$ws.Activate() | Out-Null
$ws.Range("D710").Select() | Out-Null
